# Generate Report for Handback
#
# A new handback (for d9928d69-d185-4751-8d86-e3abaa5e7c74) was processed for
# both the "zh-cn" and "de-de" target-language sheets, but it was found to be
# based on a stale version of the source file, so the row is flagged with a
# "not the latest" error instead of being marked as successfully handed back.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3e6405f88321f0bc441675298c163de0062d077/e2e/d9928d69-d185-4751-8d86-e3abaa5e7c74.md"
$errorMsg  = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d88285ae002e92593e12ea4f42a3fac7af5735f/e2e/d9928d69-d185-4751-8d86-e3abaa5e7c74.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3e6405f88321f0bc441675298c163de0062d077/e2e/d9928d69-d185-4751-8d86-e3abaa5e7c74.md."

# -----------------------------------------------------------------
# Sheet "zh-cn" - row 7 (d9928d69-d185-4751-8d86-e3abaa5e7c74)
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Latest Target File (I7) now points at the source .md, styled + linked the
# same way the other "Latest Target File" hyperlinks already are.
$wsZh.Range("I7").Value = "d9928d69-d185-4751-8d86-e3abaa5e7c74.md"
$wsZh.Range("I7").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", "d9928d69-d185-4751-8d86-e3abaa5e7c74.md") | Out-Null

# Latest Handback File (J7)
$wsZh.Range("J7").Value = "d9928d69-d185-4751-8d86-e3abaa5e7c74.4b6bcf6dd1e01d594d9929fd338a0bf674619338.zh-cn.xlf"

# Latest Handback DateTime (K7)
$wsZh.Range("K7").Value = "2016-09-04 01:00:47"

# Error Detail (P7)
$wsZh.Range("P7").Value = $errorMsg

# -----------------------------------------------------------------
# Sheet "de-de" - row 7 (d9928d69-d185-4751-8d86-e3abaa5e7c74)
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Latest Handoff File (G7) gets a fresh generated xliff name
$wsDe.Range("G7").Value = "d9928d69-d185-4751-8d86-e3abaa5e7c74.4b6bcf6dd1e01d594d9929fd338a0bf674619338.de-de.xlf"

# Latest Target File (I7)
$wsDe.Range("I7").Value = "d9928d69-d185-4751-8d86-e3abaa5e7c74.md"
$wsDe.Range("I7").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", "d9928d69-d185-4751-8d86-e3abaa5e7c74.md") | Out-Null

# Latest Handback File (J7)
$wsDe.Range("J7").Value = "d9928d69-d185-4751-8d86-e3abaa5e7c74.4b6bcf6dd1e01d594d9929fd338a0bf674619338.de-de.xlf"

# Latest Handback DateTime (K7)
$wsDe.Range("K7").Value = "2016-09-04 01:00:54"

# Error Detail (P7)
$wsDe.Range("P7").Value = $errorMsg
